# "merged in negative test 1"
#
# Before:
#   ... Program Exits<bookmark _GoBack/> ]
#   [Heading1] Negative Test 1 - Invalid User Input
#   <empty paragraph>
#   [Heading1] Negative Test 2 - Invalid File Name
#
# After:
#   ... Program Exits ]                              (bookmark removed from here)
#   [Heading1] Negative Test 1 - Invalid User Input
#   Program shows user main menu
#   User inputs "J"
#   Invalid input exception thrown
#   <bookmark _GoBack/> Program shows user main menu  (bookmark relocated here)
#   <empty paragraph>
#   [Heading1] Negative Test 2 - Invalid File Name

$d = $word.ActiveDocument

# Find the "Negative Test 1 - Invalid User Input" heading paragraph and the
# (pre-existing, untouched) empty paragraph that immediately follows it.
$heading = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $txt = $d.Paragraphs.Item($i).Range.Text
    if ($txt -like "Negative Test 1*Invalid User Input*") {
        $heading = $i
        break
    }
}

# The bookmark sits at the end of the paragraph right before the heading
# ("Program Exits"); relocate it rather than leaving it behind.
$bm = $d.Bookmarks.Item("_GoBack")
$bm.Delete()

# anchorIndex is the paragraph right after the heading - originally the
# empty paragraph that separates "Negative Test 1" from "Negative Test 2".
$anchorIndex = $heading + 1

# --- insert "Program shows user main menu" ---
$anchor = $d.Paragraphs.Item($anchorIndex).Range
$insPoint = $anchor.Duplicate
$insPoint.Collapse(1)  # wdCollapseStart
$insPoint.InsertParagraphBefore()
$d.Paragraphs.Item($anchorIndex).Range.Text = "Program shows user main menu"
$anchorIndex = $anchorIndex + 1

# --- insert 'User inputs "J"' ---
$anchor = $d.Paragraphs.Item($anchorIndex).Range
$insPoint = $anchor.Duplicate
$insPoint.Collapse(1)
$insPoint.InsertParagraphBefore()
$d.Paragraphs.Item($anchorIndex).Range.Text = "User inputs “J”"
$anchorIndex = $anchorIndex + 1

# --- insert "Invalid input exception thrown" ---
$anchor = $d.Paragraphs.Item($anchorIndex).Range
$insPoint = $anchor.Duplicate
$insPoint.Collapse(1)
$insPoint.InsertParagraphBefore()
$d.Paragraphs.Item($anchorIndex).Range.Text = "Invalid input exception thrown"
$anchorIndex = $anchorIndex + 1

# --- insert the bookmark paragraph: "Program shows user main menu" ---
$anchor = $d.Paragraphs.Item($anchorIndex).Range
$insPoint = $anchor.Duplicate
$insPoint.Collapse(1)
$insPoint.InsertParagraphBefore()
$d.Paragraphs.Item($anchorIndex).Range.Text = "Program shows user main menu"

$bmPara = $d.Paragraphs.Item($anchorIndex).Range
$bmRange = $d.Range($bmPara.Start, $bmPara.Start)
$d.Bookmarks.Add("_GoBack", $bmRange)

Write-Output "done"
